# Carga de libros pdf para leer
# Adds a new attendance-date column "E" to the "Asistencia" sheet (a new
# date, 2020-06-16 / serial 43998) and marks attendance ("X") for the
# students that attended that day, mirroring the layout/formatting of the
# existing "D" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asistencia")

# --- Copy formatting from column D so the new column reuses the same
#     (already-existing) cell styles instead of creating new ones. ---
$ws.Range("D3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null        # xlPasteFormats

$ws.Range("D4:D36").Copy() | Out-Null
$ws.Range("E4:E36").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

$ws.Range("D37").Copy() | Out-Null
$ws.Range("E37").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$excel.CutCopyMode = 0

# --- Header cell: new attendance date ---
$ws.Range("E3").Value = 43998

# --- Row data: mark attendance with "X" for present rows, leave blank otherwise ---
$presentRows = @(4,5,6,8,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,28,29,31,32,33,34,36)
$blankRows = @(7,9,26,27,30,35)

foreach ($r in $presentRows) {
    $ws.Cells.Item($r, 5).Value = "X"
}

foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 5).Value = $null
}

# --- Totals row: count of attendance marks in the new column ---
$ws.Range("E37").Formula = "=COUNTA(E4:E36)"

# --- Column width to match the other date columns ---
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth

# --- View state like the source edit (selection) ---
$ws.Range("E33").Select()
